# "fixed for some error configure file"
#
# Scene.xlsx: the CanClone flag for the DemoSummer scene (row 6) was
# mistakenly left at 0; flip it to 1.  Also leaves the selection where the
# author ended up after making the fix.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# CanClone (column M) for row 6 (DemoSummer): 0 -> 1
$ws.Range("M6").Value = 1

# Author's cursor ended up on N12 when the workbook was saved.
$ws.Range("N12").Select()
